# "textos de la memoria, falta"
# Rename the existing sheet, add a new "memorias" sheet after it, fill it
# with the review notes/growth table, and drop a note into "que poner"!E3.

$wb = $excel.ActiveWorkbook

# --- sheet1: rename Hoja1 -> "que poner" --------------------------------
$quePoner = $wb.Worksheets.Item(1)
$quePoner.Name = "que poner"

# --- add the new "memorias" sheet right after it ------------------------
$memorias = $wb.Worksheets.Add($null, $quePoner)
$memorias.Name = "memorias"

# --- header row: years 2005-2020 in B3:Q3 --------------------------------
$col = 2
foreach ($year in 2005..2020) {
    $memorias.Cells.Item(3, $col).Value = $year
    $col++
}

# --- row 4: Crecimiento + growth percentages -----------------------------
$memorias.Range("A4").Value = "Crecimiento"
$memorias.Range("F4").Value = 0.034
$memorias.Range("G4").Value = 0.041
$memorias.Range("I4").Value = 0.052
$memorias.Range("J4").Value = 0.068
$memorias.Range("K4").Value = 0.054
$memorias.Range("L4").Value = 0.048
$memorias.Range("M4").Value = 0.043
$memorias.Range("F4:G4").NumberFormat = "0.00%"
$memorias.Range("I4:M4").NumberFormat = "0.00%"

# --- row 6 ---------------------------------------------------------------
$memorias.Range("A6").Value = "Crecimiento en relación a Sudamerica"
$memorias.Range("K6").Value = "la mas alta"
$memorias.Range("L6").Value = "la mas alta"

# --- row 7 ---------------------------------------------------------------
$memorias.Range("A7").Value = "Contexto externo"
$memorias.Range("K7").Value = "Negativo"
$memorias.Range("L7").Value = "Negativo"

# --- row 8 ---------------------------------------------------------------
$memorias.Range("A8").Value = "Sectores de mayor crecimiento"
$memorias.Range("M8").Value = "Financieros, Construcción, Manufacturas"

# --- row 9 ---------------------------------------------------------------
$memorias.Range("A9").Value = "Impulsos de lado de la oferta"
$memorias.Range("L9").Value = "Sectores no extractivos"
$memorias.Range("M9").Value = "Demanda Interna, impulsos fiscales y monetarios"

# --- row 10 ----------------------------------------------------------------
$memorias.Range("A10").Value = "del lado del gasto"
$memorias.Range("L10").Value = "impulso Fiscal y Monetario"

# --- row 11 ----------------------------------------------------------------
$memorias.Range("A11").Value = "Por el lado de la demanda"
$memorias.Range("L11").Value = "Demanda interna"
$memorias.Range("M11").Value = "Consumo privado e Infraestructura"

# --- row 12 ----------------------------------------------------------------
$memorias.Range("A12").Value = "lado Fiscal"
$memorias.Range("L12").Value = "Mayor Inversión Pública"

# --- sheet1 "que poner"!E3: a note, center aligned, non-bold ---------------
$quePoner.Range("E3").Value = "Tal vez poner PII para decir quiénes nos financian"
$quePoner.Range("E3").Font.Bold = $false

# --- back to "memorias": the rest of row 9/11, and row 5 -------------------
$memorias.Range("J9").Value = "Todos los sectores"
$memorias.Range("K9").Value = "Todos los sectores"

$memorias.Range("F11").Value = "Demanda Interna"
$memorias.Range("H11").Value = "Demanda Interna"
$memorias.Range("I11").Value = "Demanda Interna"
$memorias.Range("J11").Value = "Demanda Interna"
$memorias.Range("K11").Value = "Mayor Inversión Pública"

$memorias.Range("I9").Value = "Financieros destaca"

# --- row 5 -------------------------------------------------------------
$memorias.Range("A5").Value = "Algo especial"
$memorias.Range("H5:M5").NumberFormat = "0.00%"
$memorias.Range("H5").Value = "Modelo de desarrollo basdo en la intervención estatal"
$memorias.Range("G5").Value = "Recuperación de la economía mundial"

# --- row 8 continued -----------------------------------------------------
$memorias.Range("G8").Value = "Actividades Extractivas"

# --- row 13, 14, 15 --------------------------------------------------------
$memorias.Range("A13").Value = "Detalles"
$memorias.Range("G13").Value = "La niña afectó al sector agropecuario"
$memorias.Range("G14").Value = "San Cristobal no impulso tanto"
$memorias.Range("G15").Value = "Hidrocarburos por ventas a Brasil"

# --- row 8 final cell ------------------------------------------------------
$memorias.Range("F8").Value = "destaca sectores no extractivos, y minería"

# --- cosmetics: column width + frozen header pane + selections -------------
$memorias.Columns.Item(1).ColumnWidth = 37.140625

$quePoner.Range("E4").Select()

$memorias.Activate()
$memorias.Range("B4").Select()
$excel.ActiveWindow.FreezePanes = $true
$memorias.Range("I21").Select()

$wb.Save()
